# Houses.xlsx — add a house-count "limit by" column (STR_limitBy) between
# STR_output and STR_desc, and populate it with the building that limits
# each house type's count (townHall / lumbermill / mill / stoneMason /
# foundry). The existing STR_desc column (E) shifts right to F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing STR_desc column (E) one column to the right (F),
# freeing up column E for the new STR_limitBy data.
$ws.Range("F1").Value = $ws.Range("E1").Value()
$ws.Range("F2").Value = $ws.Range("E2").Value()
$ws.Range("F3").Value = $ws.Range("E3").Value()
$ws.Range("F4").Value = $ws.Range("E4").Value()
$ws.Range("F5").Value = $ws.Range("E5").Value()
$ws.Range("F6").Value = $ws.Range("E6").Value()

# New column E: STR_limitBy — the building that caps how many of this
# house type may be built.
$ws.Range("E1").Value = "STR_limitBy"
$ws.Range("E2").Value = "townHall"
$ws.Range("E3").Value = "lumbermill"
$ws.Range("E4").Value = "mill"
$ws.Range("E5").Value = "stoneMason"
$ws.Range("E6").Value = "foundry"

# Match the author's final selection/cursor position.
$ws.Range("A6").Select()
